# This script re-purposes the "GO Struct Bits" / "SQL Bits" siena data-array
# builder workbook from the "Account" (ACCT) struct to a "Calender" (CAL)
# struct, shrinking the field list from 28 fields down to 6 fields
# (ID, Source, Date, Time, ShortName, Description), and updates the
# selection / active-sheet state left behind by the author.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("GO Struct Bits")
$ws2 = $wb.Worksheets.Item("SQL Bits")

# --- Update the short code (B2) and struct name (B4) -----------------------
$ws1.Range("B2").Value = "CAL"
$ws1.Range("B4").Value = "Calender"

# --- Replace the field list in B6:B33 ---------------------------------------
# New, shorter list of field names replacing the old "Account" fields.
$fields = @("ID", "Source", "Date", "Time", "ShortName", "Description")

for ($i = 0; $i -lt $fields.Length; $i++) {
    $row = 6 + $i
    $ws1.Cells.Item($row, 2).Value = $fields[$i]
}

# Clear out the remaining old field names (rows 12-33) so the struct only
# has the 6 new fields.
for ($row = 12; $row -le 33; $row++) {
    $ws1.Cells.Item($row, 2).Value = ""
}

# --- Update sheet selections -------------------------------------------------
# Select GO Struct Bits first and set its lingering selection to L6:L11,
# then make SQL Bits the active / tab-selected sheet with C15:H15 selected.
$ws1.Select()
$ws1.Range("L6:L11").Select()

$ws2.Select()
$ws2.Range("C15:H15").Select()
